# Auto-generated Excel COM-interop script to update Gilgamesh Profits market-data sheets
# per scheduled runner refresh (chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 7927.857
$ws.Cells.Item(64, 10).Value = 9061.625
$ws.Cells.Item(64, 12).Value = 9061.625
$ws.Cells.Item(64, 14).Value = -9557.625
# Row 67
$ws.Cells.Item(67, 8).Value = 7927.857
$ws.Cells.Item(67, 10).Value = 9061.625
$ws.Cells.Item(67, 12).Value = 9061.625
$ws.Cells.Item(67, 14).Value = -10777.625
# Row 69
$ws.Cells.Item(69, 8).Value = 8000
$ws.Cells.Item(69, 10).Value = 8000
$ws.Cells.Item(69, 12).Value = 24000
$ws.Cells.Item(69, 14).Value = -25748
# Row 70
$ws.Cells.Item(70, 8).Value = 2719.8
$ws.Cells.Item(70, 9).Value = 2666.3333
$ws.Cells.Item(70, 11).Value = 7998.999899999999
$ws.Cells.Item(70, 13).Value = -7728.999899999999
# Row 72
$ws.Cells.Item(72, 8).Value = 8000
$ws.Cells.Item(72, 10).Value = 8000
$ws.Cells.Item(72, 12).Value = 72000
$ws.Cells.Item(72, 14).Value = -80736
# Row 73
$ws.Cells.Item(73, 8).Value = 2719.8
$ws.Cells.Item(73, 9).Value = 2666.3333
$ws.Cells.Item(73, 11).Value = 7998.999899999999
$ws.Cells.Item(73, 13).Value = -7062.999899999999
# Row 80
$ws.Cells.Item(80, 8).Value = 566
$ws.Cells.Item(80, 9).Value = 582
$ws.Cells.Item(80, 10).Value = 550
$ws.Cells.Item(80, 11).Value = 1746
$ws.Cells.Item(80, 12).Value = 1650
$ws.Cells.Item(80, 13).Value = -748
$ws.Cells.Item(80, 14).Value = -3646
# Row 83
$ws.Cells.Item(83, 8).Value = 566
$ws.Cells.Item(83, 9).Value = 582
$ws.Cells.Item(83, 10).Value = 550
$ws.Cells.Item(83, 11).Value = 5238
$ws.Cells.Item(83, 12).Value = 4950
$ws.Cells.Item(83, 13).Value = -246
$ws.Cells.Item(83, 14).Value = -14934
# Row 111
$ws.Cells.Item(111, 8).Value = 923.0714
$ws.Cells.Item(111, 9).Value = 965.1111
$ws.Cells.Item(111, 10).Value = 847.4
$ws.Cells.Item(111, 11).Value = 2895.3333
$ws.Cells.Item(111, 12).Value = 2542.2
$ws.Cells.Item(111, 13).Value = 171.6667000000002
$ws.Cells.Item(111, 14).Value = -8676.200000000001
# Row 125
$ws.Cells.Item(125, 8).Value = 965
$ws.Cells.Item(125, 9).Value = 500
$ws.Cells.Item(125, 11).Value = 4500
$ws.Cells.Item(125, 13).Value = -2040
# Row 137
$ws.Cells.Item(137, 8).Value = 24215.49
$ws.Cells.Item(137, 9).Value = 27742.783
$ws.Cells.Item(137, 11).Value = 83228.349
$ws.Cells.Item(137, 13).Value = -80678.349

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 5433.5884
$ws.Cells.Item(61, 9).Value = 4305.1816
$ws.Cells.Item(61, 10).Value = 7502.3335
$ws.Cells.Item(61, 11).Value = 4305.1816
$ws.Cells.Item(61, 12).Value = 7502.3335
$ws.Cells.Item(61, 13).Value = -4093.1816
$ws.Cells.Item(61, 14).Value = -7926.3335
# Row 74
$ws.Cells.Item(74, 8).Value = 170431.16
$ws.Cells.Item(74, 9).Value = 265583.44
$ws.Cells.Item(74, 10).Value = 3914.6667
$ws.Cells.Item(74, 11).Value = 265583.44
$ws.Cells.Item(74, 12).Value = 3914.6667
$ws.Cells.Item(74, 13).Value = -264709.44
$ws.Cells.Item(74, 14).Value = -5662.6667
# Row 77
$ws.Cells.Item(77, 8).Value = 170431.16
$ws.Cells.Item(77, 9).Value = 265583.44
$ws.Cells.Item(77, 10).Value = 3914.6667
$ws.Cells.Item(77, 11).Value = 1327917.2
$ws.Cells.Item(77, 12).Value = 19573.3335
$ws.Cells.Item(77, 13).Value = -1323549.2
$ws.Cells.Item(77, 14).Value = -28309.3335
# Row 132
$ws.Cells.Item(132, 8).Value = 2959.516
$ws.Cells.Item(132, 9).Value = 1652
$ws.Cells.Item(132, 10).Value = 6718.625
$ws.Cells.Item(132, 11).Value = 4956
$ws.Cells.Item(132, 12).Value = 20155.875
$ws.Cells.Item(132, 13).Value = -2426
$ws.Cells.Item(132, 14).Value = -25215.875
# Row 134
$ws.Cells.Item(134, 8).Value = 103999
$ws.Cells.Item(134, 10).Value = 103999
$ws.Cells.Item(134, 12).Value = 103999
$ws.Cells.Item(134, 14).Value = -114139
# Row 136
$ws.Cells.Item(136, 8).Value = 5433.5884
$ws.Cells.Item(136, 9).Value = 4305.1816
$ws.Cells.Item(136, 10).Value = 7502.3335
$ws.Cells.Item(136, 11).Value = 12915.5448
$ws.Cells.Item(136, 12).Value = 22507.0005
$ws.Cells.Item(136, 13).Value = -10365.5448
$ws.Cells.Item(136, 14).Value = -27607.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Cells.Item(105, 8).Value = 14447969
$ws.Cells.Item(105, 9).Value = 1431199.8
$ws.Cells.Item(105, 11).Value = 1431199.8
$ws.Cells.Item(105, 13).Value = -1429452.8
# Row 107
$ws.Cells.Item(107, 8).Value = 1365.4
$ws.Cells.Item(107, 9).Value = 1365.4
$ws.Cells.Item(107, 11).Value = 1365.4
$ws.Cells.Item(107, 13).Value = 554.5999999999999
# Row 134
$ws.Cells.Item(134, 8).Value = 7332.4
$ws.Cells.Item(134, 9).Value = 7174.2104
$ws.Cells.Item(134, 10).Value = 7833.3335
$ws.Cells.Item(134, 11).Value = 21522.6312
$ws.Cells.Item(134, 12).Value = 23500.0005
$ws.Cells.Item(134, 13).Value = -18987.6312
$ws.Cells.Item(134, 14).Value = -28570.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Cells.Item(58, 8).Value = 3075
$ws.Cells.Item(58, 9).Value = 2485.5833
$ws.Cells.Item(58, 11).Value = 2485.5833
$ws.Cells.Item(58, 13).Value = -2282.5833
# Row 107
$ws.Cells.Item(107, 8).Value = 689.6
$ws.Cells.Item(107, 9).Value = 737
$ws.Cells.Item(107, 11).Value = 737
$ws.Cells.Item(107, 13).Value = 1183
# Row 136
$ws.Cells.Item(136, 8).Value = 3075
$ws.Cells.Item(136, 9).Value = 2485.5833
$ws.Cells.Item(136, 11).Value = 7456.749899999999
$ws.Cells.Item(136, 13).Value = -4906.749899999999

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Cells.Item(22, 8).Value = 6161.75
$ws.Cells.Item(22, 10).Value = 6161.75
$ws.Cells.Item(22, 12).Value = 18485.25
$ws.Cells.Item(22, 14).Value = -18823.25
# Row 27
$ws.Cells.Item(27, 8).Value = 6161.75
$ws.Cells.Item(27, 10).Value = 6161.75
$ws.Cells.Item(27, 12).Value = 18485.25
$ws.Cells.Item(27, 14).Value = -18689.25
# Row 86
$ws.Cells.Item(86, 8).Value = 304.2857
$ws.Cells.Item(86, 9).Value = 304.2857
$ws.Cells.Item(86, 11).Value = 912.8571000000001
$ws.Cells.Item(86, 13).Value = 273.1428999999999
# Row 89
$ws.Cells.Item(89, 8).Value = 304.2857
$ws.Cells.Item(89, 9).Value = 304.2857
$ws.Cells.Item(89, 11).Value = 2738.5713
$ws.Cells.Item(89, 13).Value = 3189.4287
# Row 139
$ws.Cells.Item(139, 8).Value = 929.06665
$ws.Cells.Item(139, 9).Value = 929.06665
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 2787.19995
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).Value = 2352.80005
$ws.Cells.Item(139, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 143.8
$ws.Cells.Item(2, 9).Value = 43.2
$ws.Cells.Item(2, 11).Value = 43.2
$ws.Cells.Item(2, 13).Value = 69.8

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 3979.5217
$ws.Cells.Item(40, 9).Value = 3918.9302
$ws.Cells.Item(40, 11).Value = 3918.9302
$ws.Cells.Item(40, 13).Value = -3782.9302
# Row 55
$ws.Cells.Item(55, 8).Value = 441.44446
$ws.Cells.Item(55, 10).Value = 661.5
$ws.Cells.Item(55, 12).Value = 661.5
$ws.Cells.Item(55, 14).Value = -1007.5
# Row 100
$ws.Cells.Item(100, 8).Value = 3465.125
$ws.Cells.Item(100, 9).Value = 3683.25
$ws.Cells.Item(100, 11).Value = 3683.25
$ws.Cells.Item(100, 13).Value = -3142.25
# Row 122
$ws.Cells.Item(122, 8).Value = 2450.0667
$ws.Cells.Item(122, 9).Value = 2371
$ws.Cells.Item(122, 11).Value = 7113
$ws.Cells.Item(122, 13).Value = -4663

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Cells.Item(2, 8).Value = 4616369
$ws.Cells.Item(2, 9).Value = 4616369
$ws.Cells.Item(2, 11).Value = 4616369
$ws.Cells.Item(2, 13).Value = -4616257
# Row 126
$ws.Cells.Item(126, 8).Value = 1093.4
$ws.Cells.Item(126, 9).Value = 1093.4
$ws.Cells.Item(126, 11).Value = 3280.2
$ws.Cells.Item(126, 13).Value = -810.2000000000003

